# Updates the cryptocurrency price/volume listing on Sheet1 to the
# latest scrape values (GitHub Actions run, 2023-01-14 11:18 UTC).
#
# Columns: D = Price, E = Volume(1h) change %. Values are stored as
# literal text in the source sheet (not numbers/percentages), so each
# assignment uses a leading apostrophe to force Excel to keep the
# string exactly as written (preserving formats like "2,497.29%" and
# trailing zeros like "0.007180") instead of auto-converting to a
# number or percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.11"
$ws.Range("E2").Value = "'5.71%"
$ws.Range("D3").Value = "'31.76"
$ws.Range("E3").Value = "'7.37%"
$ws.Range("D4").Value = "'5.205"
$ws.Range("E4").Value = "'2.05%"
$ws.Range("D5").Value = "'0.07346"
$ws.Range("E5").Value = "'8.59%"
$ws.Range("D6").Value = "'7.805"
$ws.Range("E6").Value = "'6.10%"
$ws.Range("D7").Value = "'3.740"
$ws.Range("E7").Value = "'8.44%"
$ws.Range("D8").Value = "'1.489"
$ws.Range("E8").Value = "'7.40%"
$ws.Range("D9").Value = "'0.9059"
$ws.Range("E9").Value = "'-0.93%"
$ws.Range("D10").Value = "'0.01677"
$ws.Range("E10").Value = "'2,497.29%"
$ws.Range("D11").Value = "'0.1684"
$ws.Range("E11").Value = "'5.26%"
$ws.Range("D12").Value = "'0.07517"
$ws.Range("E12").Value = "'8.54%"
$ws.Range("D13").Value = "'0.08003"
$ws.Range("E13").Value = "'3.80%"
$ws.Range("D14").Value = "'0.02953"
$ws.Range("E14").Value = "'0.99%"
$ws.Range("D15").Value = "'0.09917"
$ws.Range("E15").Value = "'10.41%"
$ws.Range("D16").Value = "'0.001492"
$ws.Range("E16").Value = "'-6.74%"
$ws.Range("E17").Value = "'1.27%"
$ws.Range("D18").Value = "'0.006478"
$ws.Range("E18").Value = "'5.16%"
$ws.Range("D19").Value = "'3.466"
$ws.Range("E19").Value = "'0.43%"
$ws.Range("D20").Value = "'2.230"
$ws.Range("E20").Value = "'-0.08%"
$ws.Range("D21").Value = "'0.3339"
$ws.Range("E21").Value = "'4.36%"
$ws.Range("D22").Value = "'0.1324"
$ws.Range("E22").Value = "'1.40%"
$ws.Range("D23").Value = "'4.360"
$ws.Range("E23").Value = "'6.36%"
$ws.Range("E24").Value = "'2.66%"
$ws.Range("D25").Value = "'0.001215"
$ws.Range("E25").Value = "'1.73%"
$ws.Range("D26").Value = "'0.004423"
$ws.Range("E26").Value = "'6.66%"
$ws.Range("D27").Value = "'0.0001294"
$ws.Range("E27").Value = "'8.00%"
$ws.Range("D28").Value = "'0.0001739"
$ws.Range("E28").Value = "'7.66%"
$ws.Range("D40").Value = "'0.04505"
$ws.Range("E40").Value = "'5.74%"
$ws.Range("D41").Value = "'0.007180"
$ws.Range("E41").Value = "'5.21%"
$ws.Range("D42").Value = "'0.1344"
$ws.Range("E42").Value = "'8.13%"
$ws.Range("D43").Value = "'0.002318"
$ws.Range("E43").Value = "'5.59%"
$ws.Range("D44").Value = "'0.01336"
$ws.Range("E44").Value = "'2.85%"
$ws.Range("D45").Value = "'0.00006054"
$ws.Range("E45").Value = "'6.24%"
$ws.Range("E46").Value = "'-3.45%"
$ws.Range("E47").Value = "'7.25%"
